$wb = $excel.ActiveWorkbook

# Update "展览" sheet (row 4: F=277->278, row 5: F=4116->4121)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 278
$ws1.Range("F5").Value = 4121

# Update "全部类型" sheet (row 4: F=277->278, row 5: F=4116->4121)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 278
$ws4.Range("F5").Value = 4121
